$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of pupae-emergence data (vial, treatment, time_hours, pupae) appended
# starting at row 92, for time_hours = 196 and time_hours = 212 cohorts.
$newRows = @(
    "1|conditioned|196|3",
    "1|unconditioned|196|0",
    "2|conditioned|196|5",
    "2|unconditioned|196|0",
    "3|conditioned|196|5",
    "3|unconditioned|196|0",
    "4|conditioned|196|6",
    "4|unconditioned|196|0",
    "5|conditioned|196|0",
    "5|unconditioned|196|0",
    "6|conditioned|196|2",
    "6|unconditioned|196|4",
    "7|conditioned|196|0",
    "7|unconditioned|196|3",
    "8|conditioned|196|0",
    "8|unconditioned|196|4",
    "9|conditioned|196|0",
    "9|unconditioned|196|4",
    "10|conditioned|196|1",
    "10|unconditioned|196|1",
    "11|conditioned|196|3",
    "11|unconditioned|196|2",
    "12|conditioned|196|1",
    "12|unconditioned|196|9",
    "13|conditioned|196|0",
    "13|unconditioned|196|2",
    "15|conditioned|196|1",
    "14|unconditioned|196|1",
    "15|conditioned|196|NULL",
    "15|unconditioned|196|0",
    "1|conditioned|212|6",
    "1|unconditioned|212|1",
    "2|conditioned|212|10",
    "2|unconditioned|212|0",
    "3|conditioned|212|7",
    "3|unconditioned|212|0",
    "4|conditioned|212|5",
    "4|unconditioned|212|5",
    "5|conditioned|212|0",
    "5|unconditioned|212|4",
    "6|conditioned|212|3",
    "6|unconditioned|212|8",
    "7|conditioned|212|2",
    "7|unconditioned|212|13",
    "8|conditioned|212|7",
    "8|unconditioned|212|10",
    "9|conditioned|212|0",
    "9|unconditioned|212|12",
    "10|conditioned|212|6",
    "10|unconditioned|212|10",
    "11|conditioned|212|1",
    "11|unconditioned|212|5",
    "12|conditioned|212|2",
    "12|unconditioned|212|10",
    "13|conditioned|212|0",
    "13|unconditioned|212|2",
    "15|conditioned|212|2",
    "14|unconditioned|212|1",
    "15|conditioned|212|NULL",
    "15|unconditioned|212|0"
)

$startRow = 92
$r = $startRow
foreach ($line in $newRows) {
    $fields = $line.Split("|")
    $vial = [int]$fields[0]
    $treatment = $fields[1]
    $timeHours = [int]$fields[2]
    $pupaeField = $fields[3]

    $ws.Cells.Item($r, 1).Value = $vial
    $ws.Cells.Item($r, 2).Value = $treatment
    $ws.Cells.Item($r, 3).Value = $timeHours
    if ($pupaeField -ne "NULL") {
        $ws.Cells.Item($r, 4).Value = [int]$pupaeField
    }

    $r = $r + 1
}

# Update the sheet view to reflect the scroll position/zoom/selection in
# effect after entering this batch of data (scrolled down so row 117 is at
# the top, zoomed to 93%, with D120 as the active cell).
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 117
$win.ScrollColumn = 1
$win.Zoom = 93
$ws.Range("D120").Select()
